{"js": "// The \"New papers: 2020\" section is trimmed down to just its first two\n// entries, and those two entries are overwritten with (copies of) two\n// entries that already exist earlier in the doc, under \"New papers:\n// 2020-08\": \"Narsaria, A. K., ...\" and \"Wang, J., ...\".\n//\n// In terms of body paragraphs (0-based):\n//   - index 26 (\"Dalla Tiezza, M., ...\")   -> becomes a copy of index 19\n//     (\"Narsaria, A. K., ...\").\n//   - index 27 (\"F\u00f6rster, A., ...\")        -> becomes a copy of index 24\n//     (\"Wang, J., ...\").\n//   - indices 28..end (16 paragraphs) are deleted entirely.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst srcAuthors = paragraphs.items[19]; // Narsaria, A. K., ...\nconst srcWang = paragraphs.items[24]; // Wang, J., ...\n\nconst srcAuthorsOoxml = srcAuthors.getOoxml();\nconst srcWangOoxml = srcWang.getOoxml();\nawait context.sync();\n\n// Each of these paragraphs consists of exactly 3 runs: plain author list,\n// bold title, plain journal/DOI info. Recover that run-level detail (text +\n// bold) from the captured OOXML so the rebuilt paragraphs keep the same\n// formatting.\nfunction extractRuns(ooxmlValue) {\n  const runRegex = /<w:r\\b[^>]*>([\\s\\S]*?)<\\/w:r>/g;\n  const runs = [];\n  let match;\n  while ((match = runRegex.exec(ooxmlValue)) !== null) {\n    const runXml = match[1];\n    const isBold = /<w:b\\s*\\/>|<w:b\\s+[^>]*\\/>|<w:b>/.test(runXml);\n    const textMatch = /<w:t[^>]*>([\\s\\S]*?)<\\/w:t>/.exec(runXml);\n    if (!textMatch) continue;\n    const text = textMatch[1]\n      .replace(/&lt;/g, \"<\")\n      .replace(/&gt;/g, \">\")\n      .replace(/&quot;/g, '\"')\n      .replace(/&apos;/g, \"'\")\n      .replace(/&amp;/g, \"&\");\n    runs.push({ text: text, bold: isBold });\n  }\n  return runs;\n}\n\nconst authorsRuns = extractRuns(srcAuthorsOoxml.value);\nconst wangRuns = extractRuns(srcWangOoxml.value);\n\nfunction rebuildParagraph(paragraph, runs) {\n  paragraph.clear();\n  for (const run of runs) {\n    const inserted = paragraph.insertText(run.text, \"End\");\n    inserted.font.bold = run.bold;\n  }\n}\n\nrebuildParagraph(paragraphs.items[26], authorsRuns);\nrebuildParagraph(paragraphs.items[27], wangRuns);\nawait context.sync();\n\n// Delete every paragraph from index 28 to the end of the body.\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = refreshedParagraphs.items.length - 1; i >= 28; i--) {\n  refreshedParagraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# The \"New papers: 2020\" section is trimmed down to just its first two\n# entries, and those two entries are overwritten with (copies of) two\n# entries that already exist earlier in the doc, under \"New papers:\n# 2020-08\": \"Narsaria, A. K., ...\" and \"Wang, J., ...\".\n#\n# In terms of Word's 1-indexed Paragraphs collection:\n#   - Paragraph 20 (\"Narsaria, A. K., ...\")  is the source for paragraph 27\n#     (\"Dalla Tiezza, M., ...\").\n#   - Paragraph 25 (\"Wang, J., ...\")         is the source for paragraph 28\n#     (\"F\u00f6rster, A., ...\").\n#   - Paragraphs 29 through the end of the document (16 paragraphs) are\n#     deleted entirely.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphRuns($paragraph) {\n    $xml = $paragraph.Range.WordOpenXML\n    $null = $xml -match \"<w:body>(.*)</w:body>\"\n    $body = $matches[1]\n\n    $runRegex = [regex]\"<w:r\\b[^>]*>(.*?)</w:r>\"\n    $runMatches = $runRegex.Matches($body)\n\n    $result = @()\n    foreach ($m in $runMatches) {\n        $runXml = $m.Groups[1].Value\n        $isBold = [bool]($runXml -match \"<w:b\\s*/>|<w:b\\s+[^>]*/>|<w:b>\")\n        if ($runXml -match \"<w:t[^>]*>(.*?)</w:t>\") {\n            $text = $matches[1]\n        } else {\n            $text = \"\"\n        }\n        $text = $text -replace \"&lt;\", \"<\"\n        $text = $text -replace \"&gt;\", \">\"\n        $text = $text -replace \"&quot;\", '\"'\n        $text = $text -replace \"&apos;\", \"'\"\n        $text = $text -replace \"&amp;\", \"&\"\n        $result += , [PSCustomObject]@{ Text = $text; Bold = $isBold }\n    }\n    return $result\n}\n\nfunction Set-ParagraphRuns($paragraph, $runs) {\n    $rng = $paragraph.Range\n    # Trim the trailing paragraph mark off the range so Find/Replace only\n    # ever touches the paragraph's own runs.\n    $rng.End = $rng.End - 1\n\n    foreach ($run in $runs) {\n        $rng.Find.ClearFormatting()\n        $rng.Find.Execute($run.Text, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0) | Out-Null\n        $found = $rng.Duplicate\n        $found.Text = $run.Text\n        $found.Bold = [int]$run.Bold\n    }\n}\n\n$srcAuthors = $d.Paragraphs.Item(20)   # Narsaria, A. K., ...\n$srcWang = $d.Paragraphs.Item(25)      # Wang, J., ...\n\n$authorsRuns = Get-ParagraphRuns $srcAuthors\n$wangRuns = Get-ParagraphRuns $srcWang\n\n$destAuthors = $d.Paragraphs.Item(27)  # Dalla Tiezza, M., ...\n$destWang = $d.Paragraphs.Item(28)     # F\u00f6rster, A., ...\n\n# Replace the destination paragraphs' run text in place, run by run, using\n# Find/Replace scoped to each paragraph so only that paragraph is touched.\nfunction Replace-ParagraphText($paragraph, $oldRuns, $newRuns) {\n    for ($i = 0; $i -lt $oldRuns.Count; $i++) {\n        $rng = $paragraph.Range\n        $rng.End = $rng.End - 1\n        $ok = $rng.Find.Execute($oldRuns[$i].Text, $false, $false, $false, $false, $false, $true, 1, $false, $newRuns[$i].Text, 2)\n    }\n}\n\n$oldAuthorsRuns = Get-ParagraphRuns $destAuthors\n$oldWangRuns = Get-ParagraphRuns $destWang\n\nReplace-ParagraphText $destAuthors $oldAuthorsRuns $authorsRuns\nReplace-ParagraphText $destWang $oldWangRuns $wangRuns\n\n# Delete paragraphs 29 through the end of the document (16 paragraphs).\n$startPara = $d.Paragraphs.Item(29)\n$tailRange = $d.Range($startPara.Range.Start, $d.Content.End)\n$tailRange.Delete()\n"}
